# Generate Report for Handoff
# Update "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" for the
# 70f8ed2f-8512-4534-854a-f54d534aa5c6.md row (row 6) across the Overview,
# zh-cn, and de-de sheets following a new handoff xliff generation run.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G is "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-29 14:47:40"

# zh-cn sheet: column H is "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-29 14:47:35"

# de-de sheet: column H is "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-29 14:47:40"
